$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27700
$ws.Range("E2").Value = 234
$ws.Range("F2").Value = 234
$ws.Range("G2").Value = 93
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 72
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 8566
$ws.Range("L2").Value = 6049
$ws.Range("M2").Value = 2517
$ws.Range("N2").Value = 2506
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 563
$ws.Range("Q2").Value = -1016
$ws.Range("R2").Value = -331
$ws.Range("S2").Value = 1030
$ws.Range("T2").Value = 278
$ws.Range("U2").Value = -1294
$ws.Range("V2").Value = 3311
$ws.Range("W2").Value = 0.84
$ws.Range("X2").Value = 0.25
$ws.Range("Y2").Value = 2.87
$ws.Range("Z2").Value = 0.89
$ws.Range("AA2").Value = 240.32
$ws.Range("AB2").Value = 377.25
$ws.Range("AC2").Value = 218
$ws.Range("AD2").Value = 27.94
$ws.Range("AE2").Value = 7572
$ws.Range("AF2").Value = 0.8
$ws.Range("AG2").Value = 85
$ws.Range("AH2").Value = 1.4
$ws.Range("AI2").Value = 39
$ws.Range("AJ2").Value = 33134683

$ws.Range("D3").Value = 22620
$ws.Range("E3").Value = 289
$ws.Range("F3").Value = 289
$ws.Range("G3").Value = -451
$ws.Range("H3").Value = -467
$ws.Range("I3").Value = -467
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 8007
$ws.Range("L3").Value = 5980
$ws.Range("M3").Value = 2027
$ws.Range("N3").Value = 2018
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 563
$ws.Range("Q3").Value = -14
$ws.Range("R3").Value = 78
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 96
$ws.Range("U3").Value = -111
$ws.Range("V3").Value = 3445
$ws.Range("W3").Value = 1.28
$ws.Range("X3").Value = -2.07
$ws.Range("Y3").Value = -20.64
$ws.Range("Z3").Value = -5.64
$ws.Range("AA3").Value = 294.98
$ws.Range("AB3").Value = 288.77
$ws.Range("AC3").Value = -1409
$ws.Range("AD3").Value = -2.77
$ws.Range("AE3").Value = 6098
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 33134683

$ws.Range("D4").Value = 25538
$ws.Range("E4").Value = 364
$ws.Range("F4").Value = 364
$ws.Range("G4").Value = 186
$ws.Range("H4").Value = 161
$ws.Range("I4").Value = 156
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 13264
$ws.Range("L4").Value = 9722
$ws.Range("M4").Value = 3543
$ws.Range("N4").Value = 3126
$ws.Range("O4").Value = 417
$ws.Range("P4").Value = 2063
$ws.Range("Q4").Value = 734
$ws.Range("R4").Value = -1582
$ws.Range("S4").Value = 802
$ws.Range("T4").Value = 46
$ws.Range("U4").Value = 687
$ws.Range("V4").Value = 5068
$ws.Range("W4").Value = 1.43
$ws.Range("X4").Value = 0.63
$ws.Range("Y4").Value = 6.07
$ws.Range("Z4").Value = 1.51
$ws.Range("AA4").Value = 274.42
$ws.Range("AB4").Value = 57.58
$ws.Range("AC4").Value = 301
$ws.Range("AD4").Value = 10.29
$ws.Range("AE4").Value = 3789
$ws.Range("AF4").Value = 0.82
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 82533764

$ws.Range("D5").Value = 33874
$ws.Range("E5").Value = 480
$ws.Range("F5").Value = 480
$ws.Range("G5").Value = 214
$ws.Range("H5").Value = 225
$ws.Range("I5").Value = 212
$ws.Range("J5").Value = 13
$ws.Range("K5").Value = 13979
$ws.Range("L5").Value = 10102
$ws.Range("M5").Value = 3877
$ws.Range("N5").Value = 3674
$ws.Range("O5").Value = 202
$ws.Range("P5").Value = 2063
$ws.Range("Q5").Value = 818
$ws.Range("R5").Value = -1314
$ws.Range("S5").Value = 775
$ws.Range("T5").Value = 89
$ws.Range("U5").Value = 730
$ws.Range("V5").Value = 5990
$ws.Range("W5").Value = 1.42
$ws.Range("X5").Value = 0.66
$ws.Range("Y5").Value = 6.22
$ws.Range("Z5").Value = 1.65
$ws.Range("AA5").Value = 260.58
$ws.Range("AB5").Value = 82.70999999999999
$ws.Range("AC5").Value = 256
$ws.Range("AD5").Value = 10.9
$ws.Range("AE5").Value = 4454
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 0.89
$ws.Range("AI5").Value = 9.75
$ws.Range("AJ5").Value = 82533764

$ws.Range("D6").Value = 40585
$ws.Range("E6").Value = 561
$ws.Range("F6").Value = 561
$ws.Range("G6").Value = 304
$ws.Range("H6").Value = 211
$ws.Range("I6").Value = 205
$ws.Range("K6").Value = 13360
$ws.Range("L6").Value = 9329
$ws.Range("M6").Value = 4031
$ws.Range("N6").Value = 3826
$ws.Range("P6").Value = 2063
$ws.Range("Q6").Value = -103
$ws.Range("R6").Value = 841
$ws.Range("S6").Value = -890
$ws.Range("T6").Value = 68
$ws.Range("U6").Value = -171
$ws.Range("V6").Value = 5142
$ws.Range("W6").Value = 1.38
$ws.Range("X6").Value = 0.52
$ws.Range("Y6").Value = 5.48
$ws.Range("Z6").Value = 1.54
$ws.Range("AA6").Value = 231.44
$ws.Range("AB6").Value = 91.13
$ws.Range("AC6").Value = 249
$ws.Range("AD6").Value = 9.9
$ws.Range("AE6").Value = 4638
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 25
$ws.Range("AH6").Value = 1.01
$ws.Range("AI6").Value = 10.04
$ws.Range("AJ6").Value = 82533764

# Clear rows 7-9 data beyond A:C (diff removes all numeric data in these rows)
$ws.Range("D7:AJ9").ClearContents()
